$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C (header "Förändrad") holds a date serial for every data row (2..398).
# The diff shows this value changing uniformly from 45186 to 45188 for all rows.
$lastRow = 398
$ws.Range("C2:C$lastRow").Value = 45188
